$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the microphone component (row 8) with the new INMP441 module
$ws.Range("C8").Value = "INMP441 Microphone Module"
$ws.Range("D8").Value = "Omnidirectional MEMS mic (I²S output, 64dB SNR, 3.3V)"
$ws.Range("E8").Value = "INMP441"

# Widen column C to fit the new, longer component name
$ws.Columns.Item(3).ColumnWidth = 24

# Update the active selection to E8
$ws.Range("E8").Select()
